$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205262/1743535585915/aNjCWRFZTU-camera-video-segment-driver-1743535590915.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMPP63KZP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T170033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHAaCXVzLXdlc3QtMiJGMEQCIGLLbfeAYZjUCX7UfxFP7erVR%2Bgh9LJyey%2BJ2tDqguobAiAZdk7A1DvI4BT4CQnIbxasMS4dEO1txnZiXMtcVF%2BDPCrmAwjY%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMTIT0KjqTrwrXXS4hKroDkKmmEIhSuCy5ZggFwV9CjpVpLpWsHnp%2FQ4awVmmdxNhTNQ9tESlCchSw35Qf8YTkWRbvaFybNy8S9Za1OT%2BM70MnY9WC%2FLfg58ESjnd6mQUwdheex5p1jNhtTYvDWJc33msc%2Ff0XfmZ6cjfmpGy0ZDt49oAwWh5QC8QizCwHxZ600QqGX2qJwUxDaEQGp7dvL84oYwyEU0pkLof4RDIiQ5hNEo20tRp4QiJLADuZS1O2NIl%2BlLbbXOhbmzA5qUVsxOJTLtFQqHasRRqf21rEGkH%2BU5yeRZ7Rpm8IuYhSb3Hlc1P640BBlrQ6iff3Y9gVOHGTElSOuOquxAeGSFUwXfPwG6exUh1CcyNC8XfL3jgOvLBXaHDFais4HSkbOp%2BtAKTxNUuCbH09fPEyJ5J3gUTWnkdIipcWaUUbLEJH13PFAy6ys8b9QyWWS%2FAigu8ifi6AqYv8zRl0L0%2Bm798L0YuWWF4QFVtfmAXLaMMKUA8b48Uj4%2FV5OD2D5Kc9mEHJDT3EKCH8CAQCu6nVxpM7bCKzH2essc0FEoUrdTFaN6ELIq%2FYzDZAUw4zMVPJbwEHWWzQGW84X8Hj3jDXq7W%2FBjqmAZ2Dk5vQh%2BBkyEdF1pCnEajXGwznHfSSefswstw6c%2F1pG5P%2BN1kdBpIhPF1tQoD2y9Av0JrDsr0H%2BVVYHS2FKnyVsswR4jAmTrccoJpbKQerkfd7McGHW5jK1NLhsoMI6Kr53Pr3iqM1NqyuLb8NG6eME0BlMEjtkWUdWWUu5mrOzULTQiAnDYNCrIHRTPQrIlTXKe7v%2FYnO%2BNoDcyCeOI6Qu%2Bk0BoU%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2001%3A00%3A33%20GMT&X-Amz-Signature=0d3d3691d92c6a81bd07a80821baf4492ddccc627094ff1d80effbabfa92bf43'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535585915/3RxOLQEM5t-camera-video-segment-1743535590915.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMPP63KZP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T170033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHAaCXVzLXdlc3QtMiJGMEQCIGLLbfeAYZjUCX7UfxFP7erVR%2Bgh9LJyey%2BJ2tDqguobAiAZdk7A1DvI4BT4CQnIbxasMS4dEO1txnZiXMtcVF%2BDPCrmAwjY%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMTIT0KjqTrwrXXS4hKroDkKmmEIhSuCy5ZggFwV9CjpVpLpWsHnp%2FQ4awVmmdxNhTNQ9tESlCchSw35Qf8YTkWRbvaFybNy8S9Za1OT%2BM70MnY9WC%2FLfg58ESjnd6mQUwdheex5p1jNhtTYvDWJc33msc%2Ff0XfmZ6cjfmpGy0ZDt49oAwWh5QC8QizCwHxZ600QqGX2qJwUxDaEQGp7dvL84oYwyEU0pkLof4RDIiQ5hNEo20tRp4QiJLADuZS1O2NIl%2BlLbbXOhbmzA5qUVsxOJTLtFQqHasRRqf21rEGkH%2BU5yeRZ7Rpm8IuYhSb3Hlc1P640BBlrQ6iff3Y9gVOHGTElSOuOquxAeGSFUwXfPwG6exUh1CcyNC8XfL3jgOvLBXaHDFais4HSkbOp%2BtAKTxNUuCbH09fPEyJ5J3gUTWnkdIipcWaUUbLEJH13PFAy6ys8b9QyWWS%2FAigu8ifi6AqYv8zRl0L0%2Bm798L0YuWWF4QFVtfmAXLaMMKUA8b48Uj4%2FV5OD2D5Kc9mEHJDT3EKCH8CAQCu6nVxpM7bCKzH2essc0FEoUrdTFaN6ELIq%2FYzDZAUw4zMVPJbwEHWWzQGW84X8Hj3jDXq7W%2FBjqmAZ2Dk5vQh%2BBkyEdF1pCnEajXGwznHfSSefswstw6c%2F1pG5P%2BN1kdBpIhPF1tQoD2y9Av0JrDsr0H%2BVVYHS2FKnyVsswR4jAmTrccoJpbKQerkfd7McGHW5jK1NLhsoMI6Kr53Pr3iqM1NqyuLb8NG6eME0BlMEjtkWUdWWUu5mrOzULTQiAnDYNCrIHRTPQrIlTXKe7v%2FYnO%2BNoDcyCeOI6Qu%2Bk0BoU%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2001%3A00%3A33%20GMT&X-Amz-Signature=d6e5c3de17309ea2befcbcefe3c177f0ce980161ee94b8ff9a2947d2e75e8877'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535541388/Uf5dS3fS3z-camera-video-segment-driver-1743535543888.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMPP63KZP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T170033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHAaCXVzLXdlc3QtMiJGMEQCIGLLbfeAYZjUCX7UfxFP7erVR%2Bgh9LJyey%2BJ2tDqguobAiAZdk7A1DvI4BT4CQnIbxasMS4dEO1txnZiXMtcVF%2BDPCrmAwjY%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMTIT0KjqTrwrXXS4hKroDkKmmEIhSuCy5ZggFwV9CjpVpLpWsHnp%2FQ4awVmmdxNhTNQ9tESlCchSw35Qf8YTkWRbvaFybNy8S9Za1OT%2BM70MnY9WC%2FLfg58ESjnd6mQUwdheex5p1jNhtTYvDWJc33msc%2Ff0XfmZ6cjfmpGy0ZDt49oAwWh5QC8QizCwHxZ600QqGX2qJwUxDaEQGp7dvL84oYwyEU0pkLof4RDIiQ5hNEo20tRp4QiJLADuZS1O2NIl%2BlLbbXOhbmzA5qUVsxOJTLtFQqHasRRqf21rEGkH%2BU5yeRZ7Rpm8IuYhSb3Hlc1P640BBlrQ6iff3Y9gVOHGTElSOuOquxAeGSFUwXfPwG6exUh1CcyNC8XfL3jgOvLBXaHDFais4HSkbOp%2BtAKTxNUuCbH09fPEyJ5J3gUTWnkdIipcWaUUbLEJH13PFAy6ys8b9QyWWS%2FAigu8ifi6AqYv8zRl0L0%2Bm798L0YuWWF4QFVtfmAXLaMMKUA8b48Uj4%2FV5OD2D5Kc9mEHJDT3EKCH8CAQCu6nVxpM7bCKzH2essc0FEoUrdTFaN6ELIq%2FYzDZAUw4zMVPJbwEHWWzQGW84X8Hj3jDXq7W%2FBjqmAZ2Dk5vQh%2BBkyEdF1pCnEajXGwznHfSSefswstw6c%2F1pG5P%2BN1kdBpIhPF1tQoD2y9Av0JrDsr0H%2BVVYHS2FKnyVsswR4jAmTrccoJpbKQerkfd7McGHW5jK1NLhsoMI6Kr53Pr3iqM1NqyuLb8NG6eME0BlMEjtkWUdWWUu5mrOzULTQiAnDYNCrIHRTPQrIlTXKe7v%2FYnO%2BNoDcyCeOI6Qu%2Bk0BoU%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2001%3A00%3A33%20GMT&X-Amz-Signature=953864c14e2935fad4b19ba593bf4fbdb58fed3befc24a5fe09283dd17952240'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532925726/R4KzSHnOTc-camera-video-segment-driver-1743532928226.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMPP63KZP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T170033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHAaCXVzLXdlc3QtMiJGMEQCIGLLbfeAYZjUCX7UfxFP7erVR%2Bgh9LJyey%2BJ2tDqguobAiAZdk7A1DvI4BT4CQnIbxasMS4dEO1txnZiXMtcVF%2BDPCrmAwjY%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMTIT0KjqTrwrXXS4hKroDkKmmEIhSuCy5ZggFwV9CjpVpLpWsHnp%2FQ4awVmmdxNhTNQ9tESlCchSw35Qf8YTkWRbvaFybNy8S9Za1OT%2BM70MnY9WC%2FLfg58ESjnd6mQUwdheex5p1jNhtTYvDWJc33msc%2Ff0XfmZ6cjfmpGy0ZDt49oAwWh5QC8QizCwHxZ600QqGX2qJwUxDaEQGp7dvL84oYwyEU0pkLof4RDIiQ5hNEo20tRp4QiJLADuZS1O2NIl%2BlLbbXOhbmzA5qUVsxOJTLtFQqHasRRqf21rEGkH%2BU5yeRZ7Rpm8IuYhSb3Hlc1P640BBlrQ6iff3Y9gVOHGTElSOuOquxAeGSFUwXfPwG6exUh1CcyNC8XfL3jgOvLBXaHDFais4HSkbOp%2BtAKTxNUuCbH09fPEyJ5J3gUTWnkdIipcWaUUbLEJH13PFAy6ys8b9QyWWS%2FAigu8ifi6AqYv8zRl0L0%2Bm798L0YuWWF4QFVtfmAXLaMMKUA8b48Uj4%2FV5OD2D5Kc9mEHJDT3EKCH8CAQCu6nVxpM7bCKzH2essc0FEoUrdTFaN6ELIq%2FYzDZAUw4zMVPJbwEHWWzQGW84X8Hj3jDXq7W%2FBjqmAZ2Dk5vQh%2BBkyEdF1pCnEajXGwznHfSSefswstw6c%2F1pG5P%2BN1kdBpIhPF1tQoD2y9Av0JrDsr0H%2BVVYHS2FKnyVsswR4jAmTrccoJpbKQerkfd7McGHW5jK1NLhsoMI6Kr53Pr3iqM1NqyuLb8NG6eME0BlMEjtkWUdWWUu5mrOzULTQiAnDYNCrIHRTPQrIlTXKe7v%2FYnO%2BNoDcyCeOI6Qu%2Bk0BoU%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2001%3A00%3A33%20GMT&X-Amz-Signature=da5d1c8c132cfcce4d7eb8df06ab897ec43ea907af9eabf49653d71b11c367ae'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532886131/8jJIZWs74D-camera-video-segment-driver-1743532888631.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMPP63KZP%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T170033Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHAaCXVzLXdlc3QtMiJGMEQCIGLLbfeAYZjUCX7UfxFP7erVR%2Bgh9LJyey%2BJ2tDqguobAiAZdk7A1DvI4BT4CQnIbxasMS4dEO1txnZiXMtcVF%2BDPCrmAwjY%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMTIT0KjqTrwrXXS4hKroDkKmmEIhSuCy5ZggFwV9CjpVpLpWsHnp%2FQ4awVmmdxNhTNQ9tESlCchSw35Qf8YTkWRbvaFybNy8S9Za1OT%2BM70MnY9WC%2FLfg58ESjnd6mQUwdheex5p1jNhtTYvDWJc33msc%2Ff0XfmZ6cjfmpGy0ZDt49oAwWh5QC8QizCwHxZ600QqGX2qJwUxDaEQGp7dvL84oYwyEU0pkLof4RDIiQ5hNEo20tRp4QiJLADuZS1O2NIl%2BlLbbXOhbmzA5qUVsxOJTLtFQqHasRRqf21rEGkH%2BU5yeRZ7Rpm8IuYhSb3Hlc1P640BBlrQ6iff3Y9gVOHGTElSOuOquxAeGSFUwXfPwG6exUh1CcyNC8XfL3jgOvLBXaHDFais4HSkbOp%2BtAKTxNUuCbH09fPEyJ5J3gUTWnkdIipcWaUUbLEJH13PFAy6ys8b9QyWWS%2FAigu8ifi6AqYv8zRl0L0%2Bm798L0YuWWF4QFVtfmAXLaMMKUA8b48Uj4%2FV5OD2D5Kc9mEHJDT3EKCH8CAQCu6nVxpM7bCKzH2essc0FEoUrdTFaN6ELIq%2FYzDZAUw4zMVPJbwEHWWzQGW84X8Hj3jDXq7W%2FBjqmAZ2Dk5vQh%2BBkyEdF1pCnEajXGwznHfSSefswstw6c%2F1pG5P%2BN1kdBpIhPF1tQoD2y9Av0JrDsr0H%2BVVYHS2FKnyVsswR4jAmTrccoJpbKQerkfd7McGHW5jK1NLhsoMI6Kr53Pr3iqM1NqyuLb8NG6eME0BlMEjtkWUdWWUu5mrOzULTQiAnDYNCrIHRTPQrIlTXKe7v%2FYnO%2BNoDcyCeOI6Qu%2Bk0BoU%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2001%3A00%3A33%20GMT&X-Amz-Signature=0de0b69dbc6b9245ccd4d9fa3396fb25cc40702d2b63350de53a40fb96c8d950'
